$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Equipos": update points ("Puntos al finalizar la primera rueda")
# and missing-home-games ("Localías faltantes") columns (D, E) for the new
# standings after the second batch of results was generated.
# ---------------------------------------------------------------------------
$wsEquipos = $wb.Worksheets.Item("Equipos")

$equiposUpdates = @(
    @{Row=2; D=10; E=3},
    @{Row=3; D=6},
    @{Row=4; D=11},
    @{Row=6; D=5; E=2},
    @{Row=7; D=2}
)

foreach ($u in $equiposUpdates) {
    if ($u.ContainsKey("D")) {
        $wsEquipos.Cells.Item($u.Row, 4).Value = $u.D
    }
    if ($u.ContainsKey("E")) {
        $wsEquipos.Cells.Item($u.Row, 5).Value = $u.E
    }
}

# ---------------------------------------------------------------------------
# Sheet "Resultados": new champs - regenerate Local/Visita/Resultado for
# every match row (column B "Fecha" and column A "Jornada" headers stay put).
# ---------------------------------------------------------------------------
$wsResultados = $wb.Worksheets.Item("Resultados")

$resultadosUpdates = @(
    @{Row=3;  Local="E"; Visita="A"; Resultado="2:3"},
    @{Row=4;  Local="C"; Visita="F"; Resultado="3:0"},
    @{Row=5;  Local="D"; Visita="B"; Resultado="1:5"},
    @{Row=7;  Local="F"; Visita="E"; Resultado="0:0"},
    @{Row=8;  Local="B"; Visita="A"; Resultado="2:6"},
    @{Row=9;  Local="C"; Visita="D"; Resultado="5:3"},
    @{Row=11; Local="F"; Visita="B"; Resultado="1:5"},
    @{Row=12; Local="A"; Visita="D"; Resultado="9:2"},
    @{Row=13; Local="E"; Visita="C"; Resultado="2:3"},
    @{Row=15; Local="C"; Visita="A"; Resultado="4:1"},
    @{Row=16; Local="D"; Visita="F"; Resultado="2:0"},
    @{Row=17; Local="B"; Visita="E"; Resultado="4:2"},
    @{Row=19; Local="A"; Visita="F"; Resultado="1:2"},
    @{Row=20; Local="B"; Visita="C"; Resultado="3:4"},
    @{Row=21; Local="D"; Visita="E"; Resultado="4:0"},
    @{Row=23; Local="E"; Visita="B"; Resultado="2:1"},
    @{Row=24; Local="A"; Visita="D"; Resultado="3:1"},
    @{Row=25; Local="F"; Visita="C"; Resultado="1:2"},
    @{Row=27; Local="D"; Visita="F"; Resultado="1:0"},
    @{Row=28; Local="E"; Visita="A"; Resultado="2:3"},
    @{Row=29; Local="C"; Visita="B"; Resultado="2:2"},
    @{Row=31; Local="F"; Visita="A"; Resultado="2:3"},
    @{Row=32; Local="C"; Visita="E"; Resultado="4:4"},
    @{Row=33; Local="B"; Visita="D"; Resultado="3:1"},
    @{Row=35; Local="F"; Visita="B"; Resultado="2:2"},
    @{Row=36; Local="A"; Visita="C"; Resultado="3:4"},
    @{Row=37; Local="D"; Visita="E"; Resultado="3:2"},
    @{Row=39; Local="B"; Visita="A"; Resultado="4:4"},
    @{Row=40; Local="E"; Visita="F"; Resultado="1:1"},
    @{Row=41; Local="D"; Visita="C"; Resultado="0:2"}
)

foreach ($u in $resultadosUpdates) {
    $wsResultados.Cells.Item($u.Row, 3).Value = $u.Local
    $wsResultados.Cells.Item($u.Row, 4).Value = $u.Visita
    $wsResultados.Cells.Item($u.Row, 5).Value = $u.Resultado
}
